# Insert a new weekly price record for "Terminal Hortofrutícola Agro Chillán"
# (Mango) above the current row 154, shifting the existing rows 154-174 down
# to 155-175 (matches the author's "Fruta / hortaliza, semanal" weekly
# refresh commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 154..174 down to 155..175, leaving a blank row 154 to populate.
$ws.Rows(154).Insert()

$ws.Range("A154").Value = 7
$ws.Range("B154").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C154").Value = "Ñuble"
$ws.Range("D154").Value = 45142
$ws.Range("E154").Value = 16
$ws.Range("F154").Value = "Fruta"
$ws.Range("G154").Value = 100108
$ws.Range("H154").Value = "Tropicales y subtropicales"
$ws.Range("I154").Value = 100108002
$ws.Range("J154").Value = "Mango"
$ws.Range("K154").Value = "Sin especificar"
$ws.Range("L154").Value = "Primera"
$ws.Range("M154").Value = 50
$ws.Range("N154").Value = 9000
$ws.Range("O154").Value = 9000
$ws.Range("P154").Value = 9000
$ws.Range("Q154").Value = "`$/bandeja 4 kilos"
$ws.Range("R154").Value = "Brasil"
$ws.Range("S154").Value = 2250
$ws.Range("T154").Value = 4
